# "WEIGHTED AVERAGES WORK!! wooo"
#
# - Match Data: fill in the (previously blank) Teleop/Auton columns for the
#   early rows (2-10) with the "no data yet" sentinel -1, matching the rest
#   of the sheet.
# - Drive Team Data: the weighted-average columns (wAvgtele/wAvgauto/wAvgpen)
#   finally compute real numbers instead of placeholder zeros, and a stray
#   "E+Z" label crept into J2.
# - Leave the UI sitting where the author left it: Match Data active/selected
#   at C25, Per Member Data selection at A9 (no longer the active tab),
#   Drive Team Data selection at I12.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Match Data: backfill G2:H10 (Teleop / Auton) with -1
# ---------------------------------------------------------------------
$wsMatch = $wb.Worksheets.Item("Match Data")
$wsMatch.Range("G2:H10").Value = -1

# ---------------------------------------------------------------------
# Drive Team Data: weighted averages now populated, plus the J2 label
# ---------------------------------------------------------------------
$wsDrive = $wb.Worksheets.Item("Drive Team Data")

$wsDrive.Range("E2").Value = 64.13627450980393
$wsDrive.Range("F2").Value = 36.37936507936508
$wsDrive.Range("G2").Value = 22.807936507936507
$wsDrive.Range("J2").Value = "E+Z"

$wsDrive.Range("E3").Value = 88.47222222222221
$wsDrive.Range("F3").Value = 46.875
$wsDrive.Range("G3").Value = 36.75

$wsDrive.Range("E4").Value = 75.03571428571428
$wsDrive.Range("F4").Value = 28.160714285714285
$wsDrive.Range("G4").Value = 28.42857142857143

[void]$wsDrive.Range("I12").Select()

# ---------------------------------------------------------------------
# Per Member Data: no longer the active tab, cursor left at A9
# ---------------------------------------------------------------------
$wsMember = $wb.Worksheets.Item("Per Member Data")
[void]$wsMember.Range("A9").Select()

# ---------------------------------------------------------------------
# Match Data becomes the active sheet again, cursor at C25
# ---------------------------------------------------------------------
[void]$wsMatch.Activate()
[void]$wsMatch.Range("C25").Select()
